$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1853281853281853
$ws.Range("C2").Value = 0.5598455598455598
$ws.Range("J2").Value = 0.007722007722007722
$ws.Range("P2").Value = 0.1428571428571428
$ws.Range("S2").Value = 0.1042471042471042
$ws.Range("C3").Value = 0.04666666666666667
$ws.Range("J3").Value = 0.006666666666666667
$ws.Range("P3").Value = 0.76
$ws.Range("S3").Value = 0.1866666666666667
$ws.Range("P4").Value = 0.6666666666666666
$ws.Range("S4").Value = 0.3333333333333333
$ws.Range("P5").Value = 0.8
$ws.Range("S5").Value = 0.2
$ws.Range("B6").Value = 0.05213270142180094
$ws.Range("D6").Value = 0.01895734597156398
$ws.Range("E6").Value = 0.004739336492890996
$ws.Range("F6").Value = 0.03791469194312796
$ws.Range("J6").Value = 0.2037914691943128
$ws.Range("O6").Value = 0.009478672985781991
$ws.Range("Q6").Value = 0.2417061611374408
$ws.Range("R6").Value = 0.08530805687203792
$ws.Range("S6").Value = 0.3459715639810427
$ws.Range("B7").Value = 0.09696969696969697
$ws.Range("D7").Value = 0.01212121212121212
$ws.Range("E7").Value = 0.006060606060606061
$ws.Range("F7").Value = 0.06060606060606061
$ws.Range("J7").Value = 0.07878787878787878
$ws.Range("O7").Value = 0.01212121212121212
$ws.Range("Q7").Value = 0.2303030303030303
$ws.Range("R7").Value = 0.08484848484848485
$ws.Range("S7").Value = 0.4181818181818182
$ws.Range("B8").Value = 0.0960591133004926
$ws.Range("D8").Value = 0.007389162561576354
$ws.Range("E8").Value = 0.002463054187192118
$ws.Range("F8").Value = 0.06403940886699508
$ws.Range("J8").Value = 0.05911330049261083
$ws.Range("O8").Value = 0.004926108374384237
$ws.Range("Q8").Value = 0.20935960591133
$ws.Range("R8").Value = 0.1083743842364532
$ws.Range("S8").Value = 0.4482758620689655
$ws.Range("B9").Value = 0.04854368932038835
$ws.Range("D9").Value = 0.04854368932038835
$ws.Range("F9").Value = 0.06796116504854369
$ws.Range("J9").Value = 0.05339805825242718
$ws.Range("O9").Value = 0.009708737864077669
$ws.Range("Q9").Value = 0.2330097087378641
$ws.Range("R9").Value = 0.1116504854368932
$ws.Range("S9").Value = 0.4271844660194175
$ws.Range("B10").Value = 0.1050488599348534
$ws.Range("D10").Value = 0.01628664495114007
$ws.Range("E10").Value = 0.001628664495114007
$ws.Range("F10").Value = 0.0741042345276873
$ws.Range("J10").Value = 0.08224755700325732
$ws.Range("O10").Value = 0.01140065146579805
$ws.Range("Q10").Value = 0.2133550488599349
$ws.Range("R10").Value = 0.08631921824104234
$ws.Range("S10").Value = 0.4096091205211727
$ws.Range("G11").Value = 0.126984126984127
$ws.Range("J11").Value = 0.07142857142857142
$ws.Range("K11").Value = 0.1904761904761905
$ws.Range("L11").Value = 0.5833333333333334
$ws.Range("S11").Value = 0.02777777777777778
$ws.Range("G12").Value = 0.7549668874172185
$ws.Range("J12").Value = 0.1788079470198675
$ws.Range("K12").Value = 0.01324503311258278
$ws.Range("L12").Value = 0.03973509933774835
$ws.Range("S12").Value = 0.01324503311258278
$ws.Range("G13").Value = 0.6470588235294118
$ws.Range("J13").Value = 0.2647058823529412
$ws.Range("S13").Value = 0.08823529411764706
$ws.Range("F15").Value = 0.02604166666666667
$ws.Range("H15").Value = 0.140625
$ws.Range("I15").Value = 0.109375
$ws.Range("J15").Value = 0.390625
$ws.Range("K15").Value = 0.06770833333333333
$ws.Range("M15").Value = 0.015625
$ws.Range("N15").Value = 0.005208333333333333
$ws.Range("O15").Value = 0.03645833333333334
$ws.Range("S15").Value = 0.2083333333333333
$ws.Range("F16").Value = 0.03468208092485549
$ws.Range("H16").Value = 0.1271676300578035
$ws.Range("I16").Value = 0.06358381502890173
$ws.Range("J16").Value = 0.4450867052023121
$ws.Range("K16").Value = 0.115606936416185
$ws.Range("M16").Value = 0.02890173410404624
$ws.Range("O16").Value = 0.04046242774566474
$ws.Range("S16").Value = 0.1445086705202312
$ws.Range("F17").Value = 0.02070393374741201
$ws.Range("H17").Value = 0.1739130434782609
$ws.Range("I17").Value = 0.07453416149068323
$ws.Range("J17").Value = 0.4285714285714285
$ws.Range("K17").Value = 0.08074534161490683
$ws.Range("M17").Value = 0.01449275362318841
$ws.Range("O17").Value = 0.07660455486542443
$ws.Range("S17").Value = 0.1304347826086956
$ws.Range("F18").Value = 0.0196078431372549
$ws.Range("H18").Value = 0.1372549019607843
$ws.Range("I18").Value = 0.1029411764705882
$ws.Range("J18").Value = 0.4803921568627451
$ws.Range("K18").Value = 0.07843137254901961
$ws.Range("M18").Value = 0.009803921568627451
$ws.Range("O18").Value = 0.06862745098039216
$ws.Range("S18").Value = 0.1029411764705882
$ws.Range("F19").Value = 0.009554140127388535
$ws.Range("H19").Value = 0.196656050955414
$ws.Range("I19").Value = 0.09315286624203821
$ws.Range("J19").Value = 0.428343949044586
$ws.Range("K19").Value = 0.09156050955414012
$ws.Range("M19").Value = 0.01353503184713376
$ws.Range("N19").Value = 0.0007961783439490446
$ws.Range("O19").Value = 0.0613057324840764
$ws.Range("S19").Value = 0.1050955414012739

$wb.Save()
